# Refresh the "source__2" query table with two newly fetched rows (139/Test1, 140/Test2)
# and refresh the dependent range references (defined name, table, view selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Source")

# --- Append the two new query result rows (UserID, FullName, Age, Email,
#     RegistrationDate, LastLoginDate, PurchaseTotal). The underlying table
#     column format is Text ("@"), so values are written/kept as text. ---
$ws.Range("A35").Value = "139"
$ws.Range("B35").Value = "Test1"
$ws.Range("C35").Value = "24"
$ws.Range("D35").Value = "test1@example.com"
$ws.Range("E35").Value = "29-07-2024 13:03:32"
$ws.Range("F35").Value = "29-07-2024 13:04:20"
$ws.Range("G35").Value = "100"

$ws.Range("A36").Value = "140"
$ws.Range("B36").Value = "Test2"
$ws.Range("C36").Value = "24"
$ws.Range("D36").Value = "test2@example.com"
$ws.Range("E36").Value = "2024-07-29"
$ws.Range("F36").Value = "2923-07-01"
$ws.Range("G36").Value = "100"

# --- Grow the query table ("source__2") to cover the new rows, which keeps
#     the table ref / autoFilter ref / sheet dimension all in sync. ---
$lo = $ws.ListObjects.Item("source__2")
$lo.Resize($ws.Range("A1:G36"))

# --- Update the ExternalData_1 defined name (scoped to the Source sheet) so
#     it still spans the full refreshed result set. ---
$nm = $wb.Names.Item("Source!ExternalData_1")
$nm.RefersTo = "=Source!`$A`$1:`$G`$36"

# --- Reflect the refreshed view position/selection. ---
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 22
$win.ScrollColumn = 1
$ws.Range("C46").Select()
